$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C hold plain text (coin name / URL); D holds price text that
# can look numeric, so force a text format before assigning and restore the
# default "Normal" style afterwards so no stray style index is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.786.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.616.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.614.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("E10").Value = "  +10.72%  "

$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("E15").Value = "  +3.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.094.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.762.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.615.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "365.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.29%  "

$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("E23").Value = "  -1.96%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.33"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.46%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.745.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "576.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.07%  "

$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.369"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("E48").Value = "  -7.28%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.622"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "
